$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two cells whose content actually changed.
$ws.Range("A1").Value = "test"
$ws.Range("B1").Value = "test"

# Move the active selection to L16 as recorded in the saved view state.
$ws.Range("L16").Select()
